# Update Psen1-Notch1 NATMI LR-pair sheet with new TPM-derived values.
#
# The sheet has one row per (Sending cluster, Target cluster) combination.
# Column G/H = Ligand average/total expression value, which only depends on
# the Sending cluster (column A). Column M/N = Receptor average/total
# expression value, which only depends on the Target cluster (column D).
# Columns I/J are the ligand "specificity" values: G (or H) for the row's
# sending cluster divided by the sum of G (or H) across all 5 sending
# clusters. Columns O/P are the analogous receptor specificity values using
# M/N across all 5 target clusters. Columns Q/R/S/T are simply the products
# of the corresponding ligand/receptor values (G*M, H*N, I*O, J*P).
#
# New base ligand (G, H) values per sending cluster:
$newGH = @{
    "ECs"               = @(12.52413133333333, 37.572394)
    "FAPs"              = @(23.74081066666666, 71.222432)
    "Inflammatory-Mac"  = @(41.01852933333333, 123.055588)
    "MuSCs"             = @(14.087727, 28.175454)
    "Resolving-Mac"     = @(37.49906666666667, 112.4972)
}

# New base receptor (M, N) values per target cluster:
$newMN = @{
    "ECs"               = @(38.10639333333333, 114.31918)
    "FAPs"              = @(9.149484, 27.448452)
    "Inflammatory-Mac"  = @(15.023598, 45.070794)
    "MuSCs"             = @(23.556204, 47.112408)
    "Resolving-Mac"     = @(13.62041, 40.86123000000001)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sum across all 5 clusters for the specificity denominators.
$sumG = 0.0
$sumH = 0.0
foreach ($key in $newGH.Keys) {
    $sumG += $newGH[$key][0]
    $sumH += $newGH[$key][1]
}

$sumM = 0.0
$sumN = 0.0
foreach ($key in $newMN.Keys) {
    $sumM += $newMN[$key][0]
    $sumN += $newMN[$key][1]
}

$lastRow = 26
for ($r = 2; $r -le $lastRow; $r++) {
    $sendingCluster = $ws.Cells.Item($r, 1).Value2
    $targetCluster  = $ws.Cells.Item($r, 4).Value2

    $G = $newGH[$sendingCluster][0]
    $H = $newGH[$sendingCluster][1]
    $I = $G / $sumG
    $J = $H / $sumH

    $M = $newMN[$targetCluster][0]
    $N = $newMN[$targetCluster][1]
    $O = $M / $sumM
    $P = $N / $sumN

    $Q = $G * $M
    $R = $H * $N
    $S = $I * $O
    $T = $J * $P

    $ws.Cells.Item($r, 7).Value2  = $G    # G
    $ws.Cells.Item($r, 8).Value2  = $H    # H
    $ws.Cells.Item($r, 9).Value2  = $I    # I
    $ws.Cells.Item($r, 10).Value2 = $J    # J
    $ws.Cells.Item($r, 13).Value2 = $M    # M
    $ws.Cells.Item($r, 14).Value2 = $N    # N
    $ws.Cells.Item($r, 15).Value2 = $O    # O
    $ws.Cells.Item($r, 16).Value2 = $P    # P
    $ws.Cells.Item($r, 17).Value2 = $Q    # Q
    $ws.Cells.Item($r, 18).Value2 = $R    # R
    $ws.Cells.Item($r, 19).Value2 = $S    # S
    $ws.Cells.Item($r, 20).Value2 = $T    # T
}
